$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.478.33'
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").Value = '3.688.51'
$ws.Range("E3").Value = '  -3.05%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'682.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.27%  '
$ws.Range("D6").Value = "'162.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.46%  '
$ws.Range("D7").Value = '3.688.01'
$ws.Range("E7").Value = '  -3.06%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -4.15%  '
$ws.Range("E10").Value = '  -7.63%  '
$ws.Range("D11").Value = "'7.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.17%  '
$ws.Range("D12").Value = "'0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.54%  '
$ws.Range("E13").Value = '  -4.62%  '
$ws.Range("D14").Value = "'33.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.82%  '
$ws.Range("D15").Value = '4.312.34'
$ws.Range("E15").Value = '  -3.00%  '
$ws.Range("D16").Value = '3.691.85'
$ws.Range("E16").Value = '  -3.01%  '
$ws.Range("D17").Value = '69.509.17'
$ws.Range("E17").Value = '  -1.92%  '
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("E19").Value = '  -6.45%  '
$ws.Range("D21").Value = "'484.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.59%  '
$ws.Range("D22").Value = "'9.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.668"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.34%  '
$ws.Range("D24").Value = "'80.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.36%  '
$ws.Range("D25").Value = '3.833.77'
$ws.Range("E25").Value = '  -3.02%  '
$ws.Range("E26").Value = '  -8.45%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  -4.31%  '
$ws.Range("E29").Value = '  -6.95%  '
$ws.Range("E30").Value = '  -8.00%  '
$ws.Range("D31").Value = "'2.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.79%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = "'6.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.80%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = "'2.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.20%  '
$ws.Range("E34").Value = '  -3.03%  '
$ws.Range("D35").Value = "'27.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.38%  '
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").Value = '3.658.51'
$ws.Range("E37").Value = '  -2.92%  '
$ws.Range("E38").Value = '  -5.48%  '
$ws.Range("E39").Value = '  +7.35%  '
$ws.Range("D40").Value = "'0.0940"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.57%  '
$ws.Range("D41").Value = "'2.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.52%  '
$ws.Range("D44").Value = "'0.956"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.09%  '
$ws.Range("D45").Value = "'161.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.39%  '
$ws.Range("D46").Value = "'48.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.50%  '
$ws.Range("D47").Value = "'30.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.36%  '
$ws.Range("E48").Value = '  -11.54%  '
$ws.Range("D49").Value = "'0.000289"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.44%  '
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("E51").Value = '  -5.42%  '
